$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update A3 text: "E 52nd Ave" -> "Pandora Ave"
$ws.Range("A3").Value = "On April 9th, 2021 at 2050 hrs victim was walking home. An unknown SOC male approached her from behind, groped her buttocks, and ran away westbound on Pandora Ave. Victim called police. Victim stated this was the second time this had occurred in two weeks with the previous incident not being reported. Police attended and conducted an area search with negative results. Video canvas conducted with possible sources located. Safety plan created, and victim services offered."

# 2. Add new row A11 with new report text
$ws.Range("A11").Value = "On April 19, 2021 at approximately 3:40pm, Victim Grade 11 at Killarney Secondary School was walking to a grocery store at Champlain Square near the intersection of Pandora Ave/Kerr Street, Vancouver, BC, when all of a sudden her crotch was grabbed over her shorts from behind by an unknown male, who then fled the area on foot. On April 21, 2021 after she finished her afternoon class, victim approached her SLO in the school hallway and advised him of the circumstances. Victim did not observe the male's face, but observed him to be wearing a black hoody with the hood up and black or grey sweatpants. Video canvass to be conducted and victim's unwashed shorts to be obtained and submitted for forensic processing. A safety plan was implemented and support was provided by her school counsellor in addition to the request for Victim Services at a later date."

# 3. Update selection to A7 (cosmetic, matches sheetView selection in diff)
$ws.Range("A7").Select()
